$d = $word.ActiveDocument

# --- Add _GoBack bookmark at the very start of the document (paragraph 1) ---
$startRange = $d.Range(0, 0)
$d.Bookmarks.Add("_GoBack", $startRange) | Out-Null

# --- Text content replacements (whole-paragraph Find & Replace) ---
$old = "Testing And Tailoring Cloud Storages Risk Assessment"
$new = "Testing and Tailoring Cloud Storages Risk Assessment"
$d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null

$old = "This document sets out the main risks to my MSC project and also sets out how I plan to overcome these risks if they  occour.  I have found tow main main areas of risk associated with this project. Firstly, there are risks associated with the technical aspects of the project and secondly there are risks assonated with other aspects. "
$new = "This document sets out the main risks to my MSC project and also sets out how I plan to overcome these risks if they occur. I have found two main areas of risk associated with this project. Firstly, there are risks associated with personal difficulties, and secondly there are risks associated with the technical aspects of the project."
$d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null

$old = "The first major risk that I have encountered is that  I don ont hit all of the deadlines  that I have set myself in the original project plan. I veliece the reasons why this could happen is that I have been very ambitious in thee targets that I have set myself, and it is difficult to gauge at the present time how long each piece of work is   going to take me. "
$new = "The first major risk that I have encountered is that I do not hit all of the deadlines that I have set myself in the original project plan. I believe the reasons why this could happen is that I have been very ambitious in these targets that I have set myself, and it is difficult to gauge at the present time how long each piece of work is going to take me. "
$d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null

$old = "In order to mitigate this however, I have used an Agile development methodology. This means that I can look at the progress I have made each week and check to ensure that I am still on target. I can then update my project plan and objectives based on my latest set of performance data. "
$new = "In order to mitigate this, however, I have used an Agile development methodology. This means that I can look at the progress I have made each week and check to ensure that I am still on target. I can then update my project plan and objectives based on my latest set of performance data. "
$d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null

$old = "I should also include other personal difficulties such as illness in this discussion as well. This is because if I am ill for more than a day during the project, then it will almost certainly y have an impact on whether I am able to deliver the project successfully.  I will use the same planning risks that I have identified above to mitigate these. "
$new = "I should also include other personal difficulties such as illness in this discussion as well. This is because if I am ill for more than a day during the project, then it will almost certainly have an impact on whether I am able to deliver the project successfully. I will use the same planning risks that I have identified above to mitigate these. "
$d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null

$old = "The second risk that I  need to plan for is the technical issues that I am facing in my project. These issues come in two main flavours. "
$new = "The second set of risks that I need to plan for are the technical issues that I am facing in my project. These issues come in two main flavours. "
$d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null

$old = "Firslty, I am using several new technologies as part of my project, such s the Scala programming language, the Play framework for developing web services and the Akka concurrency framework.  I have chosen hese technologies because they are particularly well suited to the job that I need to deal with. Because the technologies are new to me however,  there is a chance that problems with them, such as the time needed to learn to use them effectively, could hold up the progress of the project. "
$new = "Firstly, I am using several new technologies as part of my project, such as the Scala programming language, the Play framework for developing web services and the Akka concurrency framework.  I have chosen these technologies because they are particularly well suited to the job that I need to deal with. Because the technologies are new to me however, there is a chance that problems with them, such as the time needed to learn to use them effectively, could hold up the progress of the project. "
$d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null

$old = "In order to mitigate this risk,  I learnt both Scala and Akka during the Project research module, so that I am relatively competent in these technologies before the project starts. This means that I can hit the ground running, and also means that I am less likely to be blocked in the project because of this. I have also acquired several reference books that can be sued to go for more information should I get stuck. "
$new = "In order to mitigate this risk, I learnt both Scala and Akka during the Project Research module, so that I was relatively competent in these technologies before the project started. This means that I can hit the ground running, and also means that I am less likely to be blocked in the project because of this. I have also acquired several reference books that can be used for more information should I get stuck. "
$d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null

$old = "The second risk that I need to deal with is the technical difficulties associated with the project itself.  This is because the emulator that I am trying to build relies on incredibly complex technologies and  underestimating the complexities of these could  severely slow down progress of the project. "
$new = "A further set of risks within the technical difficulties are dealing with the theoretical issues associated with the project itself.  This is because the emulator that I am trying to build relies on incredibly complex technologies and underestimating the complexities of these could severely slow down progress of the project. "
$d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null

$old = "I am trying to mitigate this  by following  good software design principles that will allow me to get a clear picture in my head of the eventual design and behavior of the system before I actually start implementing it.   I will also go through the  major algorithms with my supervisor to ensure that I do not make any serious mistakes while implementing the project. "
$new = "I am trying to mitigate this by following good software design principles that will allow me to get a clear picture in my head of the eventual design and behavior of the system before I actually start implementing it. I will also go through the major algorithms with my supervisor to ensure that I do not make any serious mistakes while implementing the project. "
$d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null

# --- Set justified ("both") alignment on paragraphs 4 through 23 ---
for ($i = 4; $i -le 23; $i++) {
    $d.Paragraphs.Item($i).Alignment = 3
}
